$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("apple ipad air 2 gris espacial", "https://m.media-amazon.com/images/I/71Cl5A3GNBL._AC_SL1500_.jpg"),
    @("apple ipad air 2 gris espacial", "https://m.media-amazon.com/images/I/71jEngt8lFL._AC_SL1500_.jpg"),
    @("apple ipad air 2 gris espacial", "https://m.media-amazon.com/images/I/81j3pfkLuLL._AC_SL1500_.jpg"),
    @("apple ipad air 2 gris espacial", "https://m.media-amazon.com/images/I/81tKkMJ3RpL._AC_SL1500_.jpg"),
    @("apple ipad air 2 gris espacial", "https://m.media-amazon.com/images/I/81jtzj4UGwL._AC_SL1500_.jpg"),
    @("apple ipad air 2 gris espacial", "https://m.media-amazon.com/images/I/71cQcbCT2kL._AC_SL1500_.jpg"),
    @("iphone 12 pro max oro", "https://m.media-amazon.com/images/I/71fJ-gmBZtL._AC_SL1500_.jpg"),
    @("iphone 12 pro max oro", "https://m.media-amazon.com/images/I/81kQwxUVehL._AC_SL1500_.jpg"),
    @("iphone 12 pro max oro", "https://m.media-amazon.com/images/I/81-xL4AsWTL._AC_SL1500_.jpg"),
    @("iphone 12 pro max oro", "https://m.media-amazon.com/images/I/816ps14FFTL._AC_SL1500_.jpg"),
    @("iphone 12 pro max oro", "https://m.media-amazon.com/images/I/81tojw1-0hL._AC_SL1500_.jpg"),
    @("iphone 12 pro max oro", "https://m.media-amazon.com/images/I/71KHR0Ld3PL._AC_SL1500_.jpg")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
